$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 273, pushing the existing rows 273-320
# (and everything below) down to 276-323.
$ws.Rows("273:275").Insert()

# Fill in the 3 newly inserted rows with the new weekly price records
# (Fruta, Feria Lagunitas de Puerto Montt - Kiwi).

# Row 273: Especial
$ws.Range("A273").Value2 = 4
$ws.Range("B273").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C273").Value2 = "Los Lagos"
$ws.Range("D273").Value2 = 44782
$ws.Range("E273").Value2 = 10
$ws.Range("F273").Value2 = "Fruta"
$ws.Range("G273").Value2 = 100101
$ws.Range("H273").Value2 = "Berries"
$ws.Range("I273").Value2 = 100101007
$ws.Range("J273").Value2 = "Kiwi"
$ws.Range("K273").Value2 = "Hayward"
$ws.Range("L273").Value2 = "Especial"
$ws.Range("M273").Value2 = 300
$ws.Range("N273").Value2 = 17000
$ws.Range("O273").Value2 = 17000
$ws.Range("P273").Value2 = 17000
$ws.Range("Q273").Value2 = "`$/caja 15 kilos"
$ws.Range("R273").Value2 = "Región de O'Higgins"
$ws.Range("S273").Value2 = 1133
$ws.Range("T273").Value2 = 15

# Row 274: Primera
$ws.Range("A274").Value2 = 4
$ws.Range("B274").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C274").Value2 = "Los Lagos"
$ws.Range("D274").Value2 = 44782
$ws.Range("E274").Value2 = 10
$ws.Range("F274").Value2 = "Fruta"
$ws.Range("G274").Value2 = 100101
$ws.Range("H274").Value2 = "Berries"
$ws.Range("I274").Value2 = 100101007
$ws.Range("J274").Value2 = "Kiwi"
$ws.Range("K274").Value2 = "Hayward"
$ws.Range("L274").Value2 = "Primera"
$ws.Range("M274").Value2 = 300
$ws.Range("N274").Value2 = 14000
$ws.Range("O274").Value2 = 14000
$ws.Range("P274").Value2 = 14000
$ws.Range("Q274").Value2 = "`$/caja 15 kilos"
$ws.Range("R274").Value2 = "Región de O'Higgins"
$ws.Range("S274").Value2 = 933
$ws.Range("T274").Value2 = 15

# Row 275: Segunda
$ws.Range("A275").Value2 = 4
$ws.Range("B275").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C275").Value2 = "Los Lagos"
$ws.Range("D275").Value2 = 44782
$ws.Range("E275").Value2 = 10
$ws.Range("F275").Value2 = "Fruta"
$ws.Range("G275").Value2 = 100101
$ws.Range("H275").Value2 = "Berries"
$ws.Range("I275").Value2 = 100101007
$ws.Range("J275").Value2 = "Kiwi"
$ws.Range("K275").Value2 = "Hayward"
$ws.Range("L275").Value2 = "Segunda"
$ws.Range("M275").Value2 = 300
$ws.Range("N275").Value2 = 12000
$ws.Range("O275").Value2 = 12000
$ws.Range("P275").Value2 = 12000
$ws.Range("Q275").Value2 = "`$/caja 15 kilos"
$ws.Range("R275").Value2 = "Región de O'Higgins"
$ws.Range("S275").Value2 = 800
$ws.Range("T275").Value2 = 15
